$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"2.71374700488991e-10"
$ws.Range("C2").Value = [double]"6.708468553440206e-05"
$ws.Range("D2").Value = [double]"189.6080260415259"
$ws.Range("E2").Value = [double]"2797.565817734744"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = [double]"2987.173910861226"

$ws.Range("B3").Value = [double]"0.2881169905109251"
$ws.Range("C3").Value = [double]"0.3048912486333797"
$ws.Range("D3").Value = [double]"186123.597850132"
$ws.Range("E3").Value = [double]"2797.565817734744"
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = [double]"188921.7566761059"
